# Update "want to go" counts (column F) on several rows across sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 1465
$ws1.Range("F14").Value = 2952
$ws1.Range("F17").Value = 439
$ws1.Range("F24").Value = 677
$ws1.Range("F26").Value = 255
$ws1.Range("F30").Value = 330

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value = 238
$ws2.Range("F8").Value = 286

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1465
$ws4.Range("F19").Value = 2952
$ws4.Range("F22").Value = 238
$ws4.Range("F23").Value = 439
$ws4.Range("F32").Value = 286
$ws4.Range("F34").Value = 677
$ws4.Range("F39").Value = 255
$ws4.Range("F43").Value = 330
